# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Column G ("K") values are updated in place for rows 2-36.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 9
    3  = 7
    4  = 7
    5  = 4
    6  = 4
    7  = 4
    8  = 3
    9  = 7
    10 = 5
    11 = 5
    12 = 8
    13 = 2
    14 = 4
    15 = 3
    16 = 3
    17 = 3
    18 = 6
    19 = 7
    20 = 7
    21 = 3
    22 = 4
    23 = 5
    24 = 3
    25 = 5
    26 = 9
    27 = 2
    28 = 4
    29 = 3
    30 = 5
    31 = 1
    32 = 5
    33 = 1
    34 = 2
    35 = 3
    36 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
